$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I1").Value = "% of Q Drop's"
$ws.Range("I3").Value = "'0.00%"
$ws.Range("I6").Value = "'0.00%"
$ws.Range("I9").Value = "'9.09%"
$ws.Range("I12").Value = "'0.00%"
$ws.Range("I15").Value = "'0.00%"
$ws.Range("I18").Value = "'0.00%"
$ws.Range("I21").Value = "'0.00%"
$ws.Range("I24").Value = "'0.00%"
$ws.Range("I27").Value = "'0.00%"
$ws.Range("I30").Value = "'0.00%"
$ws.Range("I33").Value = "'0.00%"
$ws.Range("I36").Value = "'0.00%"
$ws.Range("I39").Value = "'0.00%"
$ws.Range("I42").Value = "'15.38%"
$ws.Range("I45").Value = "'0.00%"
